$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Quantities of Interest" sheet: insert a new column before column E
#    (a "Status" column, populated with "active"), which pushes the old
#    E..H columns (and their header comments) one slot to the right, F..I.
# ---------------------------------------------------------------------------
$wsQoi = $wb.Worksheets.Item("Quantities of Interest")

# Capture the original header-comment text (with original author prefix)
# before we touch anything, so we can re-attach the same wording one
# column further right afterwards.
$commentText = @{}
foreach ($addr in @("E1", "F1", "G1", "H1")) {
    $commentText[$addr] = $wsQoi.Range($addr).Comment.Text()
}

# Insert the new column; cell values (and column formatting) in E:H shift
# to F:I, but comments stay anchored to their original cell addresses, so
# they now need to be re-pointed one column to the right "by hand".
$wsQoi.Columns.Item(5).Insert()

# New "Status" column header + sample value.
$wsQoi.Range("E1").Value = "Status"
$wsQoi.Range("E2").Value = "active"

# Give the new column its own (narrower) width, matching the other
# "short" columns on this sheet (column D). The host's character->pixel
# quantization can't reproduce the exact 6.140625 figure Excel itself
# would store, so use the input that lands closest to it.
$wsQoi.Columns.Item(5).ColumnWidth = 5.3

# Re-point the existing comments to their new (shifted-right) cells.
# F1/G1/H1 already hold a comment object (left behind by the column
# insert above, which doesn't move comments) - update its text in place
# so the original author + bold "del:" prefix formatting carries over.
# I1 never had a comment, so it has to be created from scratch.
$wsQoi.Range("F1").Comment.Text($commentText["E1"]) | Out-Null
$wsQoi.Range("G1").Comment.Text($commentText["F1"]) | Out-Null
$wsQoi.Range("H1").Comment.Text($commentText["G1"]) | Out-Null
$wsQoi.Range("I1").AddComment($commentText["H1"]) | Out-Null

# The column-E comment object is now a leftover duplicate (its text was
# copied onto F1 above) sitting on what is now the "Status" column, which
# shouldn't carry any comment at all - remove it.
$wsQoi.Range("E1").Comment.Delete()

# Selection / active-cell bookkeeping on this sheet.
$wsQoi.Range("F6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Make "Quantities of Interest" the active tab (was "Constant
#    Parameters"), and drop the stale "tabSelected" marker on the old tab.
# ---------------------------------------------------------------------------
$wsQoi.Activate() | Out-Null

